$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.341.06"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").Value = "1.937.88"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7225"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.82%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3317"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07235"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8112"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08097"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "1.938.03"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.484"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "30.353.06"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008236"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.900"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "2.194.01"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.981"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.754"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.382"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.37%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.345"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.440"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.178"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05211"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.292"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7515"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.836"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.448"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4551"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.038"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8478"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.828"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.463"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4188"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06041"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
